# On Pilgrimage - March 1952 (DDLW #631)
#
# Converts the old "Heading1" title + bold "By Dorothy Day" byline into a
# pandoc-style title block: a "Title"-styled paragraph with the article
# title (word-by-word runs) and an "Authors"-styled paragraph with just the
# author's name (word-by-word runs, no "By " prefix, no bold). The legacy
# bookmark that wrapped the old heading paragraph is also removed.

$d = $word.ActiveDocument

$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(2)

# Sanity-check we are editing the paragraphs we expect before mutating.
if ($p1.Range.Text.Trim() -ne "On Pilgrimage - March 1952" -or `
    $p2.Range.Text.Trim() -ne "By Dorothy Day") {
    throw "Unexpected document content; aborting."
}

# Merge the old title + byline paragraphs into simple placeholder text
# first. Replacing the *whole* combined range (rather than deleting piece
# by piece) is what actually drops the stray <w:bookmarkEnd> that sat right
# after the byline paragraph - a plain .Delete() on sub-ranges leaves it
# behind.
$full = $d.Range($p1.Range.Start, $p2.Range.End)
$full.Text = "x" + [char]13 + "x"

# The matching <w:bookmarkStart> sits at the very start of the document
# body (position 0) and is untouched by the text replacement above because
# it's a zero-width marker right at that boundary. A zero-length delete
# exactly at position 0 removes it cleanly without touching any content.
$d.Range(0, 0).Delete()

# Now rebuild paragraph 1 as the pandoc-style Title paragraph, with each
# word (and the interstitial spaces) as its own run, matching the target
# markup exactly.
$p1 = $d.Paragraphs(1)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">On</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">-</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">March</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">1952</w:t></w:r>' + `
  '</w:p>'
$p1.Range.InsertXML($titleXml)

# Rebuild paragraph 2 as the pandoc-style Authors paragraph.
$p2 = $d.Paragraphs(2)
$authorsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + `
  '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' + `
  '</w:p>'
$p2.Range.InsertXML($authorsXml)
